$d = $word.ActiveDocument

# 1. Update the placeholder ID text in the first paragraph.
$d.Content.Find.Execute("**ID__AFFARS_5302_topic_3__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5302_101__ID**", 2)

# 2. Remove the trailing space-only run that follows the ID text.
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$spaceRange = $d.Range($r.End - 2, $r.End - 1)
$spaceRange.Delete()

# 3. Adjust the paragraph's left indent (120 twips -> 225 twips == 6pt -> 11.25pt).
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# 4. Add the paragraph border (space=5 on all sides, no visible line).
$b = $p1.Range.ParagraphFormat.Borders
$b.DistanceFromTop = 5
$b.DistanceFromBottom = 5
$b.DistanceFromLeft = 5
$b.DistanceFromRight = 5
